$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-43) holds the "Förändrad" (last changed) date, stored as
# Excel serial date 45836 (2025-06-28). This automatic update bumps it by
# one day to 45837 (2025-06-29) for every row.
for ($row = 2; $row -le 43; $row++) {
    $ws.Cells.Item($row, 3).Value = 45837
}
